$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (existing B "dbExcel"/Neo4jData shifts to C,
# existing C "WebExcel"/WebData shifts to D), creating space for the new "StatQuery" column.
$ws.Columns("B").Insert()

# Match the width used by column A for the newly inserted column B.
$ws.Columns("B").ColumnWidth = 75.81640625

# Header row.
$ws.Range("B1").Value = "StatQuery"

# Data row (new stat/count query), using the same wrap-text style as A2.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Parathyroid cancer, NOS']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Update the active selection to A2.
$null = $ws.Range("A2").Select()
